# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts from the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, shared by both sheets.
$updates = @{
    3  = 119
    4  = 1630
    5  = 625
    6  = 1098
    7  = 19
    8  = 11551
    11 = 453
    12 = 366
    14 = 805
    15 = 12397
    16 = 13123
    18 = 144
    21 = 232
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
